# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (blue/orange "Office" palette)
#   ppt/theme/theme2.xml -> "Integral"     (pink/violet "Red Violet" palette)
# and the slide master / presentation-level theme relationship currently
# resolves to the "Integral" (Red Violet) palette.
#
# The authored change swaps the two palettes: the presentation's effective
# theme becomes the "Office Theme" palette (while the other part keeps the
# "Integral" colors). Re-apply that by rewriting the 12 theme colors (the
# only thing that actually differs between the two theme parts - font
# scheme and format scheme are identical) on the presentation's theme
# color scheme, from the current "Integral" values to the "Office Theme"
# values.

function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB color values are packed as 0x00BBGGRR.
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme, in MsoThemeColorSchemeIndex order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToComRgb $officeThemeHex[$i - 1]
}

# Best-effort: the authored diff also renames the theme ("Integral" ->
# "Office Theme") and its color scheme ("Red Violet" -> "Office"). Attempt
# this too in case the host surfaces it.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
